$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(744).Insert()

$ws.Cells.Item(744,1).Value = 3
$ws.Cells.Item(744,2).Value = "Femacal de La Calera"
$ws.Cells.Item(744,3).Value = "Coquimbo"
$ws.Cells.Item(744,4).Value = 45132
$ws.Cells.Item(744,5).Value = 5
$ws.Cells.Item(744,6).Value = 100112032
$ws.Cells.Item(744,7).Value = "Zapallo italiano"
$ws.Cells.Item(744,8).Value = "Sin especificar"
$ws.Cells.Item(744,9).Value = "Primera"
$ws.Cells.Item(744,10).Value = 95
$ws.Cells.Item(744,11).Value = 13500
$ws.Cells.Item(744,12).Value = 14000
$ws.Cells.Item(744,13).Value = 13763
$ws.Cells.Item(744,14).Value = "$/caja 60 unidades"
$ws.Cells.Item(744,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(744,16).Value = 229
$ws.Cells.Item(744,17).Value = 60
$ws.Cells.Item(744,18).Value = "Hortaliza"
